# Automatische test-sync: 2025-08-05 19:40:50
# Adds a new test-mail log row (#9) to the "Logs" sheet and updates the
# corresponding category tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Append the new row to the "Logs" sheet -----------------------------
$logs = $wb.Worksheets.Item("Logs")

$newRow = 50
$logs.Cells.Item($newRow, 1).Value = "Hebben jullie toevallig al iets gehoord?"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Testmail #9: Hebben jullie toevallig al iets gehoord?"
$logs.Cells.Item($newRow, 4).Value = "Overig"
$logs.Cells.Item($newRow, 5).Value = "Bedankt, we hebben dit doorgestuurd naar support@bedrijf.nl."
$logs.Cells.Item($newRow, 6).Value = "2025-08-05 19:40:09"
$logs.Cells.Item($newRow, 7).Value = "Ja"
$logs.Cells.Item($newRow, 8).Value = "Ja"
$logs.Cells.Item($newRow, 9).Value = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# --- 2. Extend the conditional-formatting ranges so they keep covering the
#        whole data range (D/G/H/I/J 2:49 -> 2:50), without touching the
#        underlying rules themselves. -----------------------------------
$ranges = @("D2:D49", "G2:G49", "H2:H49", "I2:I49", "J2:J49")
foreach ($addr in $ranges) {
    $col = $addr.Substring(0, 1)
    $newAddr = "$($col)2:$($col)50"
    $fcs = $logs.Range($addr).FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($logs.Range($newAddr))
    }
}

# --- 3. Update the "Dashboard" tally for the "Overig" category -------------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(3, 2).Value = 8
